# Hortaliza, Macroferia Regional de Talca - Choclo
# Inserts a new weekly record as row 134 (pushing the former rows
# 134-179 down to 135-180), matching the "Fruta / hortaliza, semanal"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 134..179 down to 135..180, leaving a blank row 134.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new weekly entry.
$ws.Cells.Item(134, 1).Value = 5
$ws.Cells.Item(134, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(134, 3).Value = "Maule"
$ws.Cells.Item(134, 4).Value = 44588
$ws.Cells.Item(134, 5).Value = 7
$ws.Cells.Item(134, 6).Value = 100112024
$ws.Cells.Item(134, 7).Value = "Choclo"
$ws.Cells.Item(134, 8).Value = "Choclero"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 60000
$ws.Cells.Item(134, 11).Value = 120
$ws.Cells.Item(134, 12).Value = 120
$ws.Cells.Item(134, 13).Value = 120
$ws.Cells.Item(134, 14).Value = "`$/unidad"
$ws.Cells.Item(134, 15).Value = "Región del Maule"
$ws.Cells.Item(134, 16).Value = 120
$ws.Cells.Item(134, 17).Value = 1
$ws.Cells.Item(134, 18).Value = "Hortaliza"
